$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name) to reflect the new "through" date
$ws.Name = "Through 2021-09-13"

# Update the label in column A for the September row to reflect the new date
$ws.Range("A10").Value = "September (through 09-13)"

# Update the September row (row 10) values for each year column (B..H)
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = 23
$ws.Range("D10").Value = 32
$ws.Range("E10").Value = 26
$ws.Range("F10").Value = 30
$ws.Range("G10").Value = 46
$ws.Range("H10").Value = 66

# Update the Total row (row 11) values for each year column (B..H)
$ws.Range("B11").Value = 206
$ws.Range("C11").Value = 404
$ws.Range("D11").Value = 583
$ws.Range("E11").Value = 516
$ws.Range("F11").Value = 379
$ws.Range("G11").Value = 830
$ws.Range("H11").Value = 1137
